$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -11.9732
$ws.Range("B7").Value = 5.641399999999999
$ws.Range("D7").Value = -7.590299999999993
$ws.Range("A9").Value = -21.85180000000001
$ws.Range("D10").Value = -8.259399999999994
$ws.Range("B12").Value = 5.328799999999998
$ws.Range("A13").Value = -22.49130000000001
$ws.Range("D13").Value = -8.512099999999998
$ws.Range("B14").Value = 5.6496
$ws.Range("C15").Value = -14.5534
$ws.Range("A16").Value = -21.47509999999997
$ws.Range("D16").Value = -9.065300000000002
$ws.Range("A18").Value = -22.41330000000001
$ws.Range("B19").Value = 8.656600000000001
$ws.Range("A20").Value = -20.05439999999998
$ws.Range("D20").Value = -7.082699999999997
$ws.Range("D24").Value = -7.647699999999999
$ws.Range("A26").Value = -21.17619999999998
$ws.Range("B26").Value = 3.829300000000002
$ws.Range("A27").Value = -21.37389999999997
$ws.Range("B27").Value = 5.322600000000002
$ws.Range("C28").Value = -12.6023
$ws.Range("A29").Value = -21.71249999999999
$ws.Range("B29").Value = 5.487299999999998
$ws.Range("D32").Value = -9.161699999999987
$ws.Range("C33").Value = -11.51879999999999
$ws.Range("A35").Value = -20.4952
$ws.Range("C35").Value = -11.86210000000001
$ws.Range("A36").Value = -20.27709999999999
$ws.Range("B37").Value = 8.685700000000004
$ws.Range("B38").Value = 4.277699999999999
$ws.Range("C38").Value = -12.4513
$ws.Range("D39").Value = -7.294800000000005
$ws.Range("C43").Value = -14.02549999999999
$ws.Range("C44").Value = -13.74469999999998
$ws.Range("A45").Value = -21.8833
$ws.Range("C45").Value = -13.46769999999999
$ws.Range("B47").Value = 5.463000000000003
$ws.Range("C47").Value = -12.47869999999999
$ws.Range("D47").Value = -7.914399999999998
$ws.Range("D48").Value = -7.042099999999997
$ws.Range("B51").Value = 6.087400000000006
$ws.Range("C51").Value = -11.8088
$ws.Range("B52").Value = 5.057
$ws.Range("D52").Value = -7.895299999999999
$ws.Range("C54").Value = -13.221
$ws.Range("A55").Value = -22.5934
$ws.Range("B55").Value = 4.782699999999996
$ws.Range("D56").Value = -7.936299999999997
$ws.Range("A57").Value = -21.92010000000001
$ws.Range("C57").Value = -13.84919999999999
$ws.Range("C62").Value = -14.60470000000001
$ws.Range("C63").Value = -11.4763
$ws.Range("C67").Value = -10.65799999999999
$ws.Range("A69").Value = -21.6704
$ws.Range("B69").Value = 5.365599999999997
$ws.Range("B70").Value = 6.172300000000003
$ws.Range("C70").Value = -11.57339999999999
$ws.Range("A76").Value = -22.2933
$ws.Range("B76").Value = 5.189599999999996
$ws.Range("A78").Value = -19.84259999999998
$ws.Range("B81").Value = 5.2137
$ws.Range("C81").Value = -12.5128
$ws.Range("A82").Value = -21.99570000000001
$ws.Range("A83").Value = -22.05099999999999
$ws.Range("B83").Value = 6.098300000000003
$ws.Range("D84").Value = -8.923200000000001
$ws.Range("C88").Value = -12.3412
$ws.Range("A93").Value = -20.66019999999998
$ws.Range("B94").Value = 5.463899999999996
$ws.Range("C96").Value = -11.8934
$ws.Range("A97").Value = -21.6764
$ws.Range("C99").Value = -12.09299999999999
$ws.Range("B100").Value = 4.653099999999998
$ws.Range("D100").Value = -8.520599999999996
$ws.Range("D101").Value = -7.723399999999995
$ws.Range("B102").Value = 8.625500000000002
